$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.288.01"
$ws.Range("E2").Value = "  +4.59%  "

$ws.Range("D3").Value = "1.715.21"
$ws.Range("E3").Value = "  +3.89%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "0.9985"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "240.53"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +2.93%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.9991"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.01%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4720"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -0.99%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.2644"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +2.94%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.06235"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +2.11%  "

$ws.Range("D10").Value = "1.703.39"
$ws.Range("E10").Value = "  +3.16%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.07084"
$cell.Style = "Normal"

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "15.26"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +5.97%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "4.420"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +2.62%  "

$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.5901"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +2.45%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "76.09"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +3.29%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "0.9998"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +0.04%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.9997"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +0.10%  "

$ws.Range("D18").Value = "26.284.29"
$ws.Range("E18").Value = "  +4.63%  "

$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "11.61"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +2.70%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "0.000006790"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +2.43%  "

$ws.Range("D21").Value = "1.922.00"
$ws.Range("E21").Value = "  +3.83%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "4.572"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +5.67%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "8.845"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +4.63%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "5.344"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +1.94%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "135.28"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +0.59%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "15.19"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +1.56%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "1.404"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +2.30%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "1.759"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +7.15%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "106.45"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +2.61%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "4.026"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +3.27%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "3.688"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +4.27%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.07774"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +2.24%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.04409"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +3.24%  "

$ws.Range("E34").Value = "  +1.37%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.6208"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +4.49%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.9702"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +3.57%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.9188"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +8.20%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "111.66"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +13.19%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "2.405"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -6.96%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "1.910"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +6.79%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +0.17%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.01465"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -1.08%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.3814"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +3.59%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "5.156"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +11.13%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.1141"
$cell.Style = "Normal"

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "6.242"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +2.51%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.05298"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +1.54%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "30.70"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +5.66%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "7.695"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +6.69%  "

$ws.Range("E50").Value = "  +2.14%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.3381"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +2.89%  "
